$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Unprotect()

# Update the confidential / date notice text in A10 (change date 2021-04-09 -> 2021-04-21)
$ws.Range("A10").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-04-21 for illustrative purposes only and are subject to change."
$ws.Rows.Item(10).EntireRow.AutoFit()

# Update Weight (D) and Percent Change (E) values for rows 2-7
$ws.Range("D2").Value = 0.4873437643162154
$ws.Range("E2").Value = 0.003110419906687589

$ws.Range("D3").Value = 0.3321492881704975
$ws.Range("E3").Value = 0.009660421545667264

$ws.Range("D4").Value = 0.09696413645834827
$ws.Range("E4").Value = 0.00683945284377252

$ws.Range("D5").Value = 0.05423255112388719
$ws.Range("E5").Value = 0.003332950235604981

$ws.Range("D6").Value = 0.02931025993105155
$ws.Range("E6").Value = 0.01531100478468894

$ws.Range("D7").Value = 0.9999999999999998
$ws.Range("E7").Value = 0.006017251448684036

$ws.Protect()
